$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flip the sign of the pub_debt figures in column C (rows 2-28) so they
# are plotted as positive values.
for ($r = 2; $r -le 28; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $orig = $cell.Value2
    $cell.Value = -1 * $orig
}

# Move the active selection to H3 (scratch area used while building the
# R/Python plotting code for agenda_1 and agenda_3).
$ws.Range("H3").Select()
